$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("37:38").Insert()
$ws.Range("A37").Value = "blynk"
$ws.Range("B37").Value = "state"
$ws.Range("F37").Value = "long text"
